$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3: j = 1 coefficients
$ws.Range("B3").Value = "(0.051)"
$ws.Range("C3").Value = "(0.052)"

# Row 4: j = 2 coefficients
$ws.Range("B4").Value = "(0.054)"
$ws.Range("C4").Value = "(0.057)"

# Row 6: j = 4 coefficients
$ws.Range("B6").Value = "(0.070)"
$ws.Range("C6").Value = "(0.074)"

# Row 7: j = 5 coefficients
$ws.Range("B7").Value = "(0.079)"
$ws.Range("C7").Value = "(0.083)"

# Row 9: x block, first column
$ws.Range("B9").Value = "(0.051)"
$ws.Range("C9").Value = "(0.051)"

# Row 13: var(M1[i])
$ws.Range("B13").Value = "(0.025)"
$ws.Range("C13").Value = "(0.027)"

# Row 14: var(M2[i>id])
$ws.Range("B14").Value = "(0.092)"
$ws.Range("C14").Value = "(0.093)"

# Row 16: second x block, C column
$ws.Range("C16").Value = "(0.062)"

# Row 17: M1[i]
$ws.Range("C17").Value = "(0.155)"

# Row 18: M2[i>id]
$ws.Range("C18").Value = "(0.030)"

# Row 19: Intercept
$ws.Range("C19").Value = "(0.046)"

# Row 20: ln_p
$ws.Range("C20").Value = "(0.025)"

# Row 21: Number of observations
$ws.Range("B21").Value = 8415
$ws.Range("C21").Value = 10251
